# Apply updated crypto price/volume figures to columns D (Price) and E (Volume(1h))
# for rows 2-51, matching the scraped-data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.149.45'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '3.385.38'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''587.23'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").Value = '''180.17'
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("E9").Value = '  +6.82%  '
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("D11").Value = '''48.63'
$ws.Range("E11").Value = '  +3.75%  '
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("D13").Value = '''677.11'
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").Value = '''8.63'
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").Value = '3.929.74'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = '69.239.56'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '3.394.73'
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = '''17.67'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").Value = '''0.901'
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").Value = '''5.42'
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("D23").Value = '''17.11'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").Value = '''103.43'
$ws.Range("E24").Value = '  +4.28%  '
$ws.Range("D25").Value = '''3.92'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").Value = '''9.59'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '''34.14'
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("D29").Value = '''8.71'
$ws.Range("D30").Value = '''6.99'
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("D31").Value = '''11.19'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("E32").Value = '  +11.95%  '
$ws.Range("D33").Value = '''554.12'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").Value = '''57.96'
$ws.Range("E35").Value = '  +1.10%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '3.690.37'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").Value = '''0.140'
$ws.Range("E38").Value = '  +6.62%  '
$ws.Range("D39").Value = '''35.06'
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("E40").Value = '  +1.28%  '
$ws.Range("D41").Value = '0.0₃0702'
$ws.Range("D42").Value = '''2.67'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  +3.64%  '
$ws.Range("E45").Value = '  -1.17%  '
$ws.Range("D46").Value = '''2.66'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '''0.130'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  +5.31%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").Value = '''2.59'
$ws.Range("E51").Value = '  -0.67%  '

Write-Host "Updated 79 cells"
